$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A9:E9").NumberFormat = "@"

$ws.Range("A9").Value = "2023-12-05"
$ws.Range("B9").Value = "Final Test of rest"
$ws.Range("C9").Value = "TnjThis is the final test to see if everything resets. "
$ws.Range("D9").Value = "images/1701834132599"
$ws.Range("E9").Value = "Stressed"

$ws.Range("A9:E9").ClearFormats()
